# Update the "Förändrad" (Changed) date column (C) for every data row
# from serial date 45180 (2023-09-11) to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C301").Value = 45181
